# Weekly update: insert a new price record for "Haba" at Vega Central
# Mapocho de Santiago, pushing the existing rows 335-400 down by one and
# appending the prior last record as the new row 401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 335 (shifts 335:400 -> 336:401, dimension
# grows to A1:R401 automatically).
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row with the latest week's observation.
$ws.Cells.Item(335, 1).Value = 9
$ws.Cells.Item(335, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(335, 3).Value = "Metropolitana"
$ws.Cells.Item(335, 4).Value = 45209
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112026
$ws.Cells.Item(335, 7).Value = "Haba"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 70
$ws.Cells.Item(335, 11).Value = 11000
$ws.Cells.Item(335, 12).Value = 12000
$ws.Cells.Item(335, 13).Value = 11500
$ws.Cells.Item(335, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(335, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(335, 16).Value = 460
$ws.Cells.Item(335, 17).Value = 25
$ws.Cells.Item(335, 18).Value = "Hortaliza"
